# Add a new review row (row 10) to Sheet1, mirroring the layout/styling of
# the existing rows, add the two mailto hyperlinks for the new email
# cells, widen column B slightly, and move the selection to F10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 10 values -------------------------------------------------------
$ws.Range("A10").Value = "com.singleton.strechy"
$ws.Range("B10").Value = "stretchy"
$ws.Range("C10").Value = "dony1098765432@gmail.com"
$ws.Range("D10").Value = "sixsevensix67676@gmail.com"
$ws.Range("E10").Value = "27/5/2019 15:59"
$ws.Range("F10").Value = "Nice car game. Very addictive and fun!"

# --- Hyperlinks for the two email cells -----------------------------------
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:dony1098765432@gmail.com", "", "", "dony1098765432@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:sixsevensix67676@gmail.com", "", "", "sixsevensix67676@gmail.com")

# --- Re-apply the same formatting used by the rest of the table ----------
# (Hyperlinks.Add stamps its own "Hyperlink" look; put it back to match
# the other rows, which render as plain centered text.)
$ws.Range("A10").Font.Name = "Mangal"
$ws.Range("A10").Font.Size = 10

$ws.Range("C10").Font.Name = "Calibri"
$ws.Range("C10").Font.Size = 11
$ws.Range("C10").Font.Underline = -4142
$ws.Range("C10").Font.Color = 0
$ws.Range("C10").HorizontalAlignment = -4108

$ws.Range("D10").Font.Name = "Calibri"
$ws.Range("D10").Font.Size = 11
$ws.Range("D10").Font.Underline = -4142
$ws.Range("D10").Font.Color = 0
$ws.Range("D10").HorizontalAlignment = -4108

$ws.Range("F10").Font.Name = "Mangal"
$ws.Range("F10").Font.Size = 10

# --- Minor column width tweak ---------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 19.45

# --- Selection lands on the last edited cell ------------------------------
$ws.Range("F10").Select()
